$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - 想去人数 (F column) updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 1820
$wsExhibit.Range("F7").Value = 122
$wsExhibit.Range("F8").Value = 126
$wsExhibit.Range("F9").Value = 3062
$wsExhibit.Range("F13").Value = 587
$wsExhibit.Range("F14").Value = 504
$wsExhibit.Range("F16").Value = 354
$wsExhibit.Range("F19").Value = 1295
$wsExhibit.Range("F23").Value = 602
$wsExhibit.Range("F24").Value = 41
$wsExhibit.Range("F28").Value = 86
$wsExhibit.Range("F30").Value = 80
$wsExhibit.Range("F31").Value = 3462
$wsExhibit.Range("F33").Value = 62
$wsExhibit.Range("F34").Value = 228
$wsExhibit.Range("F36").Value = 1682

# Sheet "全部类型" (sheet4) - 想去人数 (F column) updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 1820
$wsAll.Range("F7").Value = 122
$wsAll.Range("F8").Value = 126
$wsAll.Range("F9").Value = 3062
$wsAll.Range("F13").Value = 587
$wsAll.Range("F14").Value = 504
$wsAll.Range("F17").Value = 354
$wsAll.Range("F20").Value = 1295
$wsAll.Range("F24").Value = 602
$wsAll.Range("F25").Value = 41
$wsAll.Range("F29").Value = 86
$wsAll.Range("F31").Value = 80
$wsAll.Range("F32").Value = 3462
$wsAll.Range("F35").Value = 62
$wsAll.Range("F36").Value = 228
$wsAll.Range("F38").Value = 1682
